$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.226.80"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "3.092.01"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'522.33"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "'136.33"
$ws.Range("E6").Value = "  -3.79%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "3.092.69"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "'0.458"
$ws.Range("E9").Value = "  +2.98%  "
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").Value = "'0.400"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "3.627.55"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("D17").Value = "57.333.97"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "3.093.04"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "'5.90"
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("D20").Value = "'12.47"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "'7.86"
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").Value = "'349.33"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'68.66"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").Value = "'0.499"
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  -6.29%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'7.27"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "'5.84"
$ws.Range("E32").Value = "  -8.25%  "
$ws.Range("D33").Value = "'20.90"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").Value = "'4.86"
$ws.Range("E34").Value = "  +4.92%  "
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("D36").Value = "'158.67"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "'6.02"
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("D38").Value = "'25.58"
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").Value = "'0.0658"
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").Value = "'4.05"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").Value = "'0.695"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("D44").Value = "2.397.48"
$ws.Range("E44").Value = "  +5.63%  "
$ws.Range("D45").Value = "'36.70"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "3.133.53"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").Value = "'0.0263"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").Value = "'0.952"
$ws.Range("E49").Value = "  -4.24%  "
$ws.Range("D50").Value = "'5.96"
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").Value = "'0.765"
$ws.Range("E51").Value = "  +1.34%  "
